$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.146.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.824.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6206"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07350"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.99"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07672"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.826.68"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.956"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6647"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.31"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008977"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.842"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.125.22"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.072.25"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.04"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.42"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.232"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1429"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.479"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.65"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05572"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.088"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.837"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7341"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.213.24"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01764"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.301"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9119"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.64"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.975.61"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.62"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5086"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000117"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.17%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4018"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.092"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05756"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.19%  "
